$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price/Volume columns so that values such as
# "1.00", "7.52", "0.0782" etc. are written back as literal text instead of being
# auto-converted to numbers (which would silently drop formatting/trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '60.809.17'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '3.389.84'
$ws.Range("E3").Value = '  -2.54%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '571.59'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = '141.98'
$ws.Range("E6").Value = '  -4.35%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.389.55'
$ws.Range("E8").Value = '  -2.55%  '
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").Value = '7.52'
$ws.Range("E10").Value = '  -2.79%  '
$ws.Range("E11").Value = '  -2.61%  '
$ws.Range("D12").Value = '0.395'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '3.970.43'
$ws.Range("E13").Value = '  -2.26%  '
$ws.Range("D14").Value = '28.37'
$ws.Range("E14").Value = '  +1.16%  '
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("E16").Value = '  -2.22%  '
$ws.Range("D17").Value = '3.393.74'
$ws.Range("E17").Value = '  -2.43%  '
$ws.Range("D18").Value = '60.892.97'
$ws.Range("E18").Value = '  -1.84%  '
$ws.Range("D19").Value = '6.28'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '14.12'
$ws.Range("E20").Value = '  -2.60%  '
$ws.Range("E21").Value = '  -5.97%  '
$ws.Range("D22").Value = '388.47'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("D23").Value = '0.562'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").Value = '73.51'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '0.0000118'
$ws.Range("E26").Value = '  -5.21%  '
$ws.Range("D27").Value = '3.531.59'
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("D28").Value = '0.179'
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").Value = '7.40'
$ws.Range("E30").Value = '  -5.25%  '
$ws.Range("D31").Value = '8.01'
$ws.Range("E31").Value = '  -3.16%  '
$ws.Range("E32").Value = '  -2.06%  '
$ws.Range("E33").Value = '  -7.22%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = '23.66'
$ws.Range("E35").Value = '  -2.98%  '
$ws.Range("D36").Value = '6.96'
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("D37").Value = '167.20'
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").Value = '3.421.21'
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("E39").Value = '  -4.39%  '
$ws.Range("E40").Value = '  -5.57%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '0.0782'
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '27.95'
$ws.Range("E42").Value = '  +2.45%  '
$ws.Range("D43").Value = '0.783'
$ws.Range("E43").Value = '  -3.76%  '
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").Value = '4.44'
$ws.Range("E45").Value = '  -2.14%  '
$ws.Range("D46").Value = '41.65'
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("E47").Value = '  -2.96%  '
$ws.Range("D48").Value = '2.553.38'
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("E49").Value = '  -4.40%  '
$ws.Range("D50").Value = '23.42'
$ws.Range("E50").Value = '  +0.88%  '
$ws.Range("D51").Value = '6.85'
$ws.Range("E51").Value = '  -1.52%  '

# Restore the default "Normal" style so no stray explicit style index is left
# on cells that did not have one originally.
$ws.Range("D2:E51").Style = "Normal"
